$d = $word.ActiveDocument

# Replace the four occurrences where the old "Rat" folder name is
# referenced with the new folder name
# "brainflatmapvisualizationtool-main".

$d.Content.Find.Execute(
    "Open the .zip file downloaded, which will unzip the file and create a folder called " + [char]8220 + "Rat" + [char]8221,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Open the .zip file downloaded, which will unzip the file and create a folder called " + [char]8220 + "brainflatmapvisualizationtool-main" + [char]8221,
    2
)

$d.Content.Find.Execute(
    "In the top left corner, select the " + [char]8220 + "+ New" + [char]8221 + " button, then " + [char]8220 + "File Upload" + [char]8221 + " and select the " + [char]8220 + "Rat" + [char]8221 + " folder that you just downloaded, then select " + [char]8220 + "Upload" + [char]8221 + " once more.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the top left corner, select the " + [char]8220 + "+ New" + [char]8221 + " button, then " + [char]8220 + "Folder Upload" + [char]8221 + " and select the " + [char]8220 + "brainflatmapvisualizationtool-main" + [char]8221 + " folder that you just downloaded, then select " + [char]8220 + "Upload" + [char]8221 + " once more.",
    2
)

$d.Content.Find.Execute(
    "Once uploaded, open the " + [char]8220 + "Rat" + [char]8221 + " folder and open the " + [char]8220 + "ratInputtedData.xlsx" + [char]8221 + " file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Once uploaded, open the " + [char]8220 + "brainflatmapvisualizationtool-main" + [char]8221 + " folder and open the " + [char]8220 + "ratInputtedData.xlsx" + [char]8221 + " file.",
    2
)

$d.Content.Find.Execute(
    "View your updated heatmap in the " + [char]8220 + "Rat" + [char]8221 + " folder as " + [char]8220 + "ex.svg" + [char]8221,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "View your updated heatmap in the " + [char]8220 + "brainflatmapvisualizationtool-main" + [char]8221 + " folder as " + [char]8220 + "ex.svg" + [char]8221,
    2
)
